# Progress-Report.xlsx update — "Updated Progress Report & SRS"
#
# Test Cases sheet: iteration test-case rows finished (100% complete),
# their DONE flags now compute to TRUE, and the "still needs testing"
# remarks are cleared now that testing is done.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("User Interface")
$ws2 = $wb.Worksheets.Item("Back end")
$ws3 = $wb.Worksheets.Item("Test Cases")
$ws4 = $wb.Worksheets.Item("Paperworks")

# --- Test Cases: mark rows 3, 4, 5, 17, 20, 21 as fully complete ---
$ws3.Range("D3").Value = 1
$ws3.Range("D4").Value = 1
$ws3.Range("D5").Value = 1
$ws3.Range("D17").Value = 1
$ws3.Range("D20").Value = 1
$ws3.Range("D21").Value = 1

# Remarks for those rows no longer apply now that testing is finished.
$ws3.Range("F3").Value = ""
$ws3.Range("F4").Value = ""
$ws3.Range("F20").Value = ""
$ws3.Range("F21").Value = ""

# --- Window / selection state ---
$ws1.Range("B24").Select()
$ws2.Select()
$excel.ActiveWindow.ScrollRow = 11
$ws2.Range("E21").Select()

# "Test Cases" is the sheet that ends up active/selected.
$ws3.Activate()
$ws3.Range("F21").Select()
